$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in "Actual" (F) and "Correct" (G) results for the already-recorded
#     games (rows 171-180, date: Sat, Feb 1, 2025) ---
$results171to180 = @(
    @{ Row = 171; Actual = "Brandon Wheat Kings";    Correct = 1 },
    @{ Row = 172; Actual = "Edmonton Oil Kings";      Correct = 0 },
    @{ Row = 173; Actual = "Red Deer Rebels";         Correct = 1 },
    @{ Row = 174; Actual = "Prince Albert Raiders";   Correct = 0 },
    @{ Row = 175; Actual = "Victoria Royals";         Correct = 1 },
    @{ Row = 176; Actual = "Prince George Cougars";   Correct = 0 },
    @{ Row = 177; Actual = "Spokane Chiefs";          Correct = 0 },
    @{ Row = 178; Actual = "Swift Current Broncos";   Correct = 1 },
    @{ Row = 179; Actual = "Seattle Thunderbirds";    Correct = 1 },
    @{ Row = 180; Actual = "Vancouver Giants";        Correct = 0 }
)

foreach ($r in $results171to180) {
    $ws.Cells.Item($r.Row, 6).Value = $r.Actual
    $ws.Cells.Item($r.Row, 7).Value = $r.Correct
}

# --- New games (rows 181-187) ---
$newGames = @(
    @{ Row = 181; GameID = 1021727; Date = "Sun, Feb 2, 2025"; Home = "Calgary Hitmen";         Away = "Lethbridge Hurricanes";  Prediction = "Calgary Hitmen";       Actual = "Calgary Hitmen";      Correct = 1 },
    @{ Row = 182; GameID = 1021729; Date = "Sun, Feb 2, 2025"; Home = "Saskatoon Blades";        Away = "Moose Jaw Warriors";     Prediction = "Saskatoon Blades";     Actual = "Saskatoon Blades";    Correct = 1 },
    @{ Row = 183; GameID = 1021730; Date = "Sun, Feb 2, 2025"; Home = "Vancouver Giants";        Away = "Spokane Chiefs";         Prediction = "Spokane Chiefs";       Actual = "Spokane Chiefs";      Correct = 1 },
    @{ Row = 184; GameID = 1021728; Date = "Sun, Feb 2, 2025"; Home = "Everett Silvertips";      Away = "Kamloops Blazers";       Prediction = "Everett Silvertips";   Actual = "Everett Silvertips";  Correct = 1 },
    @{ Row = 185; GameID = 1021731; Date = "Tue, Feb 4, 2025"; Home = "Moose Jaw Warriors";      Away = "Wenatchee Wild";         Prediction = "Wenatchee Wild";       Actual = $null;                 Correct = $null },
    @{ Row = 186; GameID = 1021732; Date = "Tue, Feb 4, 2025"; Home = "Tri-City Americans";      Away = "Prince George Cougars";  Prediction = "Prince George Cougars"; Actual = $null;                Correct = $null },
    @{ Row = 187; GameID = 1021733; Date = "Tue, Feb 4, 2025"; Home = "Victoria Royals";         Away = "Kelowna Rockets";        Prediction = "Victoria Royals";      Actual = $null;                 Correct = $null }
)

foreach ($g in $newGames) {
    $ws.Cells.Item($g.Row, 1).Value = $g.GameID
    $ws.Cells.Item($g.Row, 2).Value = $g.Date
    $ws.Cells.Item($g.Row, 3).Value = $g.Home
    $ws.Cells.Item($g.Row, 4).Value = $g.Away
    $ws.Cells.Item($g.Row, 5).Value = $g.Prediction
    if ($null -ne $g.Actual) {
        $ws.Cells.Item($g.Row, 6).Value = $g.Actual
    }
    if ($null -ne $g.Correct) {
        $ws.Cells.Item($g.Row, 7).Value = $g.Correct
    }
}

# --- Update view: scroll/selection moved to reflect the newly added rows ---
$win = $excel.ActiveWindow
$win.ScrollRow = 169
$win.ScrollColumn = 1
$ws.Range("F185").Select()
